# v4.2 fix GPIOs, GPIO20 is useless
#
# - "Signal Define" sheet: fill in the Signal CTL (col E) names for the
#   EPIO-mapped GPIO rows (GPIO17-GPIO22, rows 19-24). GPIO20 has no
#   real signal, so it gets a placeholder "XXX" instead of an EPIOx name.
# - Update window/selection state: "Signal Define" becomes the active /
#   selected tab (with its selection moved to the newly edited E23:E24),
#   while "Power Control" is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$wsSignal = $wb.Worksheets.Item("Signal Define")
$wsPower  = $wb.Worksheets.Item("Power Control")

# Fill in the new Signal CTL values, in the same order the original author
# typed them (this drives shared-string allocation order).
$wsSignal.Range("E22").Value = "XXX"
$wsSignal.Range("E23").Value = "EPIO5"
$wsSignal.Range("E24").Value = "EPIO4"
$wsSignal.Range("E21").Value = "EPIO3"
$wsSignal.Range("E19").Value = "EPIO1"
$wsSignal.Range("E20").Value = "EPIO2"

# "Power Control" was previously the active/selected sheet; scroll its
# window down so row 4 is the top visible row, and leave its selection
# untouched.
$wsPower.Activate()
$excel.ActiveWindow.ScrollRow = 4
$wsPower.Range("B5:G5").Select()

# Make "Signal Define" the active/selected sheet, with its selection over
# the two freshly-edited EPIO5/EPIO4 cells.
$wsSignal.Activate()
$wsSignal.Range("E23:E24").Select()
